$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8 (INTRON/EXON) - replace observations with the "does not need a function" notes
$ws.Range("B8").Value = "Does not need a function"
$ws.Range("C8").Value = "Does not need a function."

# Row 9 (CODING_DNA) - tweak the "Returns" note wording
$ws.Range("D9").Value = "The coding sequence; concatenated exons (1 string)."

$ws.Range("E8").Value = "This task can be completed by only giving the front end the full DNA sequence and an hash of exons positions/length; both retrievable with queries; all the front end needs is to know which sbstring within the main string to highligh; same we decided to do for the restriction sites task."
$ws.Range("D8").Value = "The front end can do this task with 2 lines of code by only having the full DNA sequence and an hash of exons position. We aleady have SQL fnctions retrieving those for the front end to use.  The front end just need to place tags around each exon subtring."
$ws.Range("F8").Value = "Done"

# View state tweaks observed in the diff
$ws.Application.ActiveWindow.ScrollRow = 4
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("D10").Select()
